# Apply updated transition-matrix probabilities to Sheet1.
# These values reflect re-simulated game outcomes (more games simulated,
# faster simulate-game logic, and drafted optimization logic), which
# shifted many of the Markov-chain transition probabilities slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellUpdates = @{
    "B2" = 0.1659919028340081
    "C2" = 0.631578947368421
    "J2" = 0.03238866396761134
    "P2" = 0.1093117408906883
    "S2" = 0.06072874493927125
    "B3" = 0.01197604790419162
    "C3" = 0.03592814371257485
    "J3" = 0.02395209580838323
    "P3" = 0.8143712574850299
    "S3" = 0.1137724550898204
    "J4" = 0.05
    "P4" = 0.725
    "S4" = 0.225
    "B6" = 0.08040201005025126
    "D6" = 0.01507537688442211
    "F6" = 0.09547738693467336
    "J6" = 0.2462311557788945
    "O6" = 0.03517587939698492
    "Q6" = 0.1105527638190955
    "R6" = 0.09547738693467336
    "S6" = 0.321608040201005
    "B7" = 0.09230769230769231
    "D7" = 0.04615384615384616
    "E7" = 0.007692307692307693
    "F7" = 0.04615384615384616
    "J7" = 0.1076923076923077
    "O7" = 0.007692307692307693
    "Q7" = 0.2153846153846154
    "R7" = 0.09230769230769231
    "S7" = 0.3846153846153846
    "B8" = 0.09278350515463918
    "D8" = 0.03092783505154639
    "E8" = 0.002577319587628866
    "F8" = 0.07216494845360824
    "J8" = 0.1108247422680412
    "O8" = 0.0154639175257732
    "Q8" = 0.1855670103092784
    "R8" = 0.1185567010309278
    "S8" = 0.3711340206185567
    "B9" = 0.0935672514619883
    "D9" = 0.01754385964912281
    "F9" = 0.04093567251461988
    "J9" = 0.1052631578947368
    "O9" = 0.02923976608187134
    "Q9" = 0.2105263157894737
    "R9" = 0.1286549707602339
    "S9" = 0.3742690058479532
    "B10" = 0.113314447592068
    "D10" = 0.0169971671388102
    "E10" = 0.0009442870632672333
    "F10" = 0.08120868744098206
    "J10" = 0.1076487252124646
    "O10" = 0.01038715769593957
    "Q10" = 0.1954674220963173
    "R10" = 0.09915014164305949
    "S10" = 0.3748819641170916
    "F11" = 0.004950495049504951
    "G11" = 0.1089108910891089
    "J11" = 0.1089108910891089
    "K11" = 0.2178217821782178
    "L11" = 0.5198019801980198
    "S11" = 0.0396039603960396
    "G12" = 0.7586206896551724
    "J12" = 0.1293103448275862
    "K12" = 0.008620689655172414
    "L12" = 0.0603448275862069
    "S12" = 0.04310344827586207
    "G13" = 0.7666666666666667
    "J13" = 0.2
    "S13" = 0.03333333333333333
    "F15" = 0.02808988764044944
    "H15" = 0.1966292134831461
    "I15" = 0.03932584269662921
    "J15" = 0.3707865168539326
    "K15" = 0.01685393258426966
    "N15" = 0.005617977528089887
    "O15" = 0.07303370786516854
    "S15" = 0.2696629213483146
    "F16" = 0.01612903225806452
    "H16" = 0.1881720430107527
    "I16" = 0.06989247311827956
    "J16" = 0.3709677419354839
    "K16" = 0.1075268817204301
    "M16" = 0.01612903225806452
    "N16" = 0.005376344086021506
    "O16" = 0.04838709677419355
    "S16" = 0.1774193548387097
    "F17" = 0.0215633423180593
    "H17" = 0.1725067385444744
    "I17" = 0.137466307277628
    "J17" = 0.3908355795148248
    "K17" = 0.05121293800539083
    "M17" = 0.01617250673854448
    "O17" = 0.06199460916442048
    "S17" = 0.1482479784366577
    "F18" = 0.01470588235294118
    "H18" = 0.1813725490196078
    "I18" = 0.08333333333333333
    "J18" = 0.446078431372549
    "K18" = 0.08333333333333333
    "M18" = 0.0196078431372549
    "N18" = 0.004901960784313725
    "O18" = 0.03431372549019608
    "S18" = 0.1323529411764706
    "F19" = 0.009319664492078284
    "H19" = 0.2068965517241379
    "I19" = 0.08014911463187326
    "J19" = 0.3783783783783784
    "K19" = 0.09040074557315937
    "M19" = 0.01584342963653309
    "N19" = 0.001863932898415657
    "O19" = 0.06616961789375582
    "S19" = 0.1509785647716682
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = $cellUpdates[$addr]
}
